# Auto-generated edit script applying market-price/profit updates
# from the scheduled runner diff, across sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 784.3333
$ws.Range("I19").Value = 631
$ws.Range("J19").Value = 937.6667
$ws.Range("K19").Value = 631
$ws.Range("L19").Value = 937.6667
$ws.Range("M19").Value = -456
$ws.Range("N19").Value = -1287.6667

$ws.Range("H32").Value = 71436010
$ws.Range("I32").Value = 125009576
$ws.Range("J32").Value = 4580.6665
$ws.Range("K32").Value = 125009576
$ws.Range("L32").Value = 4580.6665
$ws.Range("M32").Value = -125009250
$ws.Range("N32").Value = -5232.6665

$ws.Range("H62").Value = 151963.28
$ws.Range("I62").Value = 151963.28
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 151963.28
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -151339.28
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 151963.28
$ws.Range("I65").Value = 151963.28
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 759816.4
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -756696.4
$ws.Range("N65").ClearContents()

$ws.Range("H70").Value = 55889.1
$ws.Range("I70").Value = 170782.83
$ws.Range("J70").Value = 6648.9287
$ws.Range("K70").Value = 512348.49
$ws.Range("L70").Value = 19946.7861
$ws.Range("M70").Value = -512078.49
$ws.Range("N70").Value = -20486.7861

$ws.Range("H73").Value = 55889.1
$ws.Range("I73").Value = 170782.83
$ws.Range("J73").Value = 6648.9287
$ws.Range("K73").Value = 512348.49
$ws.Range("L73").Value = 19946.7861
$ws.Range("M73").Value = -511412.49
$ws.Range("N73").Value = -21818.7861

$ws.Range("H96").Value = 291.57144
$ws.Range("I96").Value = 241.83333
$ws.Range("J96").Value = 590
$ws.Range("K96").Value = 725.49999
$ws.Range("L96").Value = 1770
$ws.Range("M96").Value = 647.50001
$ws.Range("N96").Value = -4516

$ws.Range("H98").Value = 3293.7144
$ws.Range("I98").Value = 3173
$ws.Range("K98").Value = 3173
$ws.Range("M98").Value = -1675

$ws.Range("H122").Value = 3293.7144
$ws.Range("I122").Value = 3173
$ws.Range("K122").Value = 9519
$ws.Range("M122").Value = -7069

$ws.Range("H132").Value = 2951.3865
$ws.Range("I132").Value = 2951.3865
$ws.Range("K132").Value = 8854.1595
$ws.Range("M132").Value = -6324.1595

$ws.Range("H137").Value = 3321.1667
$ws.Range("I137").Value = 3682.8096
$ws.Range("K137").Value = 11048.4288
$ws.Range("M137").Value = -8498.4288

$ws.Range("H138").Value = 217041.98
$ws.Range("J138").Value = 242547.33
$ws.Range("L138").Value = 727641.99
$ws.Range("N138").Value = -737921.99

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2287.5334
$ws.Range("I74").Value = 2423.8462
$ws.Range("K74").Value = 2423.8462
$ws.Range("M74").Value = -1549.8462

$ws.Range("H77").Value = 2287.5334
$ws.Range("I77").Value = 2423.8462
$ws.Range("K77").Value = 12119.231
$ws.Range("M77").Value = -7751.231

$ws.Range("H122").Value = 1591.3334
$ws.Range("I122").Value = 1345.7084
$ws.Range("K122").Value = 4037.1252
$ws.Range("M122").Value = -1587.1252

$ws.Range("H132").Value = 4019.9473
$ws.Range("I132").Value = 2538.05
$ws.Range("J132").Value = 5666.5
$ws.Range("K132").Value = 7614.150000000001
$ws.Range("L132").Value = 16999.5
$ws.Range("M132").Value = -5084.150000000001
$ws.Range("N132").Value = -22059.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 98000
$ws.Range("J52").Value = 98000
$ws.Range("L52").Value = 98000
$ws.Range("N52").Value = -98526

$ws.Range("H105").Value = 6789.4287
$ws.Range("I105").Value = 5910.8
$ws.Range("K105").Value = 5910.8
$ws.Range("M105").Value = -4163.8

$ws.Range("H121").Value = 98000
$ws.Range("J121").Value = 98000
$ws.Range("L121").Value = 98000
$ws.Range("N121").Value = -101494

$ws.Range("H134").Value = 2824.35
$ws.Range("I134").Value = 2824.35
$ws.Range("K134").Value = 8473.049999999999
$ws.Range("M134").Value = -5938.049999999999

$ws.Range("H140").Value = 184420
$ws.Range("J140").Value = 283890
$ws.Range("L140").Value = 283890
$ws.Range("N140").Value = -294250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3824.8076
$ws.Range("I31").Value = 3218.111
$ws.Range("J31").Value = 4146
$ws.Range("K31").Value = 3218.111
$ws.Range("L31").Value = 4146
$ws.Range("M31").Value = -2923.111
$ws.Range("N31").Value = -4736

$ws.Range("H34").Value = 3824.8076
$ws.Range("I34").Value = 3218.111
$ws.Range("J34").Value = 4146
$ws.Range("K34").Value = 3218.111
$ws.Range("L34").Value = 4146
$ws.Range("M34").Value = -3016.111
$ws.Range("N34").Value = -4550

$ws.Range("H58").Value = 5173.294
$ws.Range("I58").Value = 4343.15
$ws.Range("K58").Value = 4343.15
$ws.Range("M58").Value = -4140.15

$ws.Range("H134").Value = 4490.8
$ws.Range("I134").Value = 4490.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13472.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -10937.4
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 5173.294
$ws.Range("I136").Value = 4343.15
$ws.Range("K136").Value = 13029.45
$ws.Range("M136").Value = -10479.45

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1804.6182
$ws.Range("J68").Value = 1929.5
$ws.Range("L68").Value = 5788.5
$ws.Range("N68").Value = -7410.5

$ws.Range("H71").Value = 1804.6182
$ws.Range("J71").Value = 1929.5
$ws.Range("L71").Value = 17365.5
$ws.Range("N71").Value = -25477.5

$ws.Range("H75").Value = 3995
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 3995
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H86").Value = 933.4286
$ws.Range("J86").Value = 933.4286
$ws.Range("L86").Value = 2800.2858
$ws.Range("N86").Value = -5172.2858

$ws.Range("H89").Value = 933.4286
$ws.Range("J89").Value = 933.4286
$ws.Range("L89").Value = 8400.857399999999
$ws.Range("N89").Value = -20256.8574

$ws.Range("H107").Value = 2085.7334
$ws.Range("J107").Value = 2085.7334
$ws.Range("L107").Value = 6257.2002
$ws.Range("N107").Value = -10097.2002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 55000
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 4588.3228
$ws.Range("I132").Value = 5436.174
$ws.Range("K132").Value = 16308.522
$ws.Range("M132").Value = -13778.522

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 131935.11
$ws.Range("I46").Value = 19500
$ws.Range("J46").Value = 145989.5
$ws.Range("K46").Value = 19500
$ws.Range("L46").Value = 145989.5
$ws.Range("M46").Value = -19312
$ws.Range("N46").Value = -146365.5

$ws.Range("H61").Value = 4208.385
$ws.Range("J61").Value = 5553.909
$ws.Range("L61").Value = 5553.909
$ws.Range("N61").Value = -5957.909

$ws.Range("H82").Value = 9420
$ws.Range("I82").Value = 13919.625
$ws.Range("K82").Value = 13919.625
$ws.Range("M82").Value = -13558.625

$ws.Range("H85").Value = 9420
$ws.Range("I85").Value = 13919.625
$ws.Range("K85").Value = 13919.625
$ws.Range("M85").Value = -12671.625

$ws.Range("H113").Value = 4208.385
$ws.Range("J113").Value = 5553.909
$ws.Range("L113").Value = 5553.909
$ws.Range("N113").Value = -9893.909

$ws.Range("H132").Value = 7409.7085
$ws.Range("I132").Value = 7402
$ws.Range("K132").Value = 22206
$ws.Range("M132").Value = -19676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 83252.2
$ws.Range("I62").Value = 116389.86
$ws.Range("K62").Value = 116389.86
$ws.Range("M62").Value = -115765.86

$ws.Range("H65").Value = 83252.2
$ws.Range("I65").Value = 116389.86
$ws.Range("K65").Value = 581949.3
$ws.Range("M65").Value = -578829.3

$ws.Range("H126").Value = 2036
$ws.Range("I126").Value = 2064
$ws.Range("J126").Value = 1949.4546
$ws.Range("K126").Value = 6192
$ws.Range("L126").Value = 5848.3638
$ws.Range("M126").Value = -3722
$ws.Range("N126").Value = -10788.3638

$ws.Range("H132").Value = 2093.6667
$ws.Range("I132").Value = 1563.8572
$ws.Range("J132").Value = 2835.4
$ws.Range("K132").Value = 4691.571599999999
$ws.Range("L132").Value = 8506.200000000001
$ws.Range("M132").Value = -2161.571599999999
$ws.Range("N132").Value = -13566.2

$ws.Range("H136").Value = 2562.2104
$ws.Range("J136").Value = 5023
$ws.Range("L136").Value = 15069
$ws.Range("N136").Value = -20169
